$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3: Adulto/Niño selector flipped from Adult(1)/Child(0) to Adult(0)/Child(1) ---
$ws.Range("G3").Value = 0
$ws.Range("H3").Value = 1

# --- Row 4: Adulto/Niño selector flipped, count changed from 2 to 0/1 ---
$ws.Range("D4").ClearContents()
$ws.Range("G4").Value = 0
$ws.Range("H4").Value = 1

# --- Row 5: drop the leftover "2" helper input, keep the rest ---
$ws.Range("D5").ClearContents()

# --- Drop the scratch "COSTON*NORMA" column N (rows 3-5) and its header/label cells ---
$ws.Range("N2").ClearContents()
$ws.Range("N3").ClearContents()
$ws.Range("N4").ClearContents()
$ws.Range("N5").ClearContents()
$ws.Range("O5").ClearContents()

# --- Drop the scratch helper row 6 (N6/P6/Q6) entirely ---
$ws.Range("N6").ClearContents()
$ws.Range("P6").ClearContents()
$ws.Range("Q6").ClearContents()

# --- Row 7: drop the now-unused N7/P7 helper cells ---
$ws.Range("N7").ClearContents()
$ws.Range("P7").ClearContents()

# --- Row 8: drop the now-unused N8 helper cell ---
$ws.Range("N8").ClearContents()

# --- Row 9: add the new F9 input, drop the now-unused N9 helper cell ---
$ws.Range("F9").Value = 4
$ws.Range("N9").ClearContents()

# --- Drop the old scratch "CALC" area entirely (rows 11-15) ---
$ws.Range("B11:Q15").ClearContents()

# --- Update the view state to match ---
$ws.Range("I5").Select()
$excel.ActiveWindow.ScrollColumn = 3
